# Update loading_percent values for the 380 kV case (Case_2_113).
# Rows 2-25 correspond to time steps 0-23; columns B,C,E,F,G,H,I,M,N,O hold
# the per-line loading percentages that were recomputed for this case.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updatedValues = @{
    "B2" = 10.36449719779026
    "C2" = 8.623358978591462
    "E2" = 12.65707304454802
    "F2" = 16.86991607391245
    "G2" = 16.77363880896904
    "H2" = 11.15400494107102
    "I2" = 15.4657048455061
    "M2" = 13.34558865085438
    "N2" = 15.87007462938115
    "O2" = 15.36939036514849
    "B3" = 9.781117537934644
    "C3" = 8.337923472744439
    "E3" = 12.53894123891221
    "F3" = 15.89584955866815
    "G3" = 16.71939826580109
    "H3" = 11.19444605686557
    "I3" = 15.56148506668162
    "M3" = 13.06469284462087
    "N3" = 15.90016538456462
    "O3" = 15.42063440703064
    "B4" = 9.402698781237763
    "C4" = 8.156388374404026
    "E4" = 12.47088295625939
    "F4" = 15.26997757108489
    "G4" = 16.69474203887264
    "H4" = 11.22139261345906
    "I4" = 15.62394337873643
    "M4" = 12.89181407556857
    "N4" = 15.92032182493226
    "O4" = 15.45650740689446
    "B5" = 9.243468403664808
    "C5" = 8.080909205742664
    "E5" = 12.44430486528213
    "F5" = 15.008197319934
    "G5" = 16.68687430486391
    "H5" = 11.23290494786572
    "I5" = 15.65031273520113
    "M5" = 12.82136366561466
    "N5" = 15.92895899691665
    "O5" = 15.47223024649135
    "B6" = 9.21672717424045
    "C6" = 8.068287599201634
    "E6" = 12.43996227279436
    "F6" = 14.96433081551589
    "G6" = 16.68569960713246
    "H6" = 11.23484864013371
    "I6" = 15.65474672458198
    "M6" = 12.80966839124466
    "N6" = 15.93041877591827
    "O6" = 15.47490757711539
    "B7" = 9.400571582625608
    "C7" = 8.155376411751805
    "E7" = 12.47051979466703
    "F7" = 15.26647399323133
    "G7" = 16.69462710167077
    "H7" = 11.22154572200593
    "I7" = 15.62429529240474
    "M7" = 12.8908638243181
    "N7" = 15.92043659422352
    "O7" = 15.45671498528059
    "B8" = 10.16757817121572
    "C8" = 8.526286817041496
    "E8" = 12.6154317904885
    "F8" = 16.5399640634477
    "G8" = 16.75314633561889
    "H8" = 11.16750961104937
    "I8" = 15.49797226196492
    "M8" = 13.24888049094879
    "N8" = 15.88010150984763
    "O8" = 15.38614167488298
    "B9" = 11.5091520536022
    "C9" = 9.20087942558599
    "E9" = 12.93361945581471
    "F9" = 19.00274580682531
    "G9" = 16.93607747533056
    "H9" = 11.07835649470376
    "I9" = 15.2792259792872
    "M9" = 13.94340300237671
    "N9" = 15.81431054464555
    "O9" = 15.282916181304
    "B10" = 12.39373034142104
    "C10" = 9.660881012520868
    "E10" = 13.18596357980648
    "F10" = 20.67494806633232
    "G10" = 17.11115656477154
    "H10" = 11.02313913719189
    "I10" = 15.13620497846846
    "M10" = 14.44355600993458
    "N10" = 15.77404726193109
    "O10" = 15.22875483619973
    "B11" = 12.77395220743509
    "C11" = 9.861814874741498
    "E11" = 13.30429501503863
    "F11" = 21.3917225636224
    "G11" = 17.19936520153037
    "H11" = 11.00025964215478
    "I11" = 15.07499041055798
    "M11" = 14.6678151469412
    "N11" = 15.75747541915428
    "O11" = 15.20886822587203
    "B12" = 12.91473028638625
    "C12" = 9.936664399025258
    "E12" = 13.34956930202903
    "F12" = 21.65686569030329
    "G12" = 17.23397060350579
    "H12" = 10.99191838698402
    "I12" = 15.05236407062119
    "M12" = 14.75218407510447
    "N12" = 15.75145024647512
    "O12" = 15.2020242187067
    "B13" = 12.88455387261889
    "C13" = 9.920599980112474
    "E13" = 13.33979871546775
    "F13" = 21.60004134736742
    "G13" = 17.22646471588208
    "H13" = 10.99370046135209
    "I13" = 15.0572123859632
    "M13" = 14.73403959660058
    "N13" = 15.75273675563719
    "O13" = 15.20346761081727
    "B14" = 12.78559837724323
    "C14" = 9.867997863664103
    "E14" = 13.30801069928009
    "F14" = 21.4136618050453
    "G14" = 17.20218828212352
    "H14" = 10.99956693030782
    "I14" = 15.0731178121343
    "M14" = 14.67476777665692
    "N14" = 15.75697471301777
    "O14" = 15.20829138643197
    "B15" = 12.72456771689219
    "C15" = 9.835614912463422
    "E15" = 13.28859882134613
    "F15" = 21.29868154950795
    "G15" = 17.18747397871338
    "H15" = 11.00320236002979
    "I15" = 15.08293256607506
    "M15" = 14.63838764187032
    "N15" = 15.75960315371075
    "O15" = 15.21133560002169
    "B16" = 12.36843431573013
    "C16" = 9.64757823660487
    "E16" = 13.17829779817789
    "F16" = 20.62722412089977
    "G16" = 17.10556187369517
    "H16" = 11.02467948578478
    "I16" = 15.14028299292646
    "M16" = 14.42882719100155
    "N16" = 15.77516531765805
    "O16" = 15.23015038122338
    "B17" = 12.14426578658191
    "C17" = 9.530059791627322
    "E17" = 13.11150694229949
    "F17" = 20.20408069597325
    "G17" = 17.05748557746532
    "H17" = 11.03842902788372
    "I17" = 15.176451592104
    "M17" = 14.29937149518354
    "N17" = 15.78515850225448
    "O17" = 15.24291213132193
    "B18" = 12.013243055735
    "C18" = 9.461685076175456
    "E18" = 13.07342640694332
    "F18" = 19.95656407809801
    "G18" = 17.03064143149555
    "H18" = 11.04654811822325
    "I18" = 15.19761677963524
    "M18" = 14.2246095589276
    "N18" = 15.79107051777764
    "O18" = 15.25069942137259
    "B19" = 11.96852311170593
    "C19" = 9.438401795025444
    "E19" = 13.06059199671625
    "F19" = 19.87204792380568
    "G19" = 17.02169207219632
    "H19" = 11.04933327137421
    "I19" = 15.20484508185094
    "M19" = 14.19924715235053
    "N19" = 15.79310044379887
    "O19" = 15.25341274149199
    "B20" = 12.16834502003559
    "C20" = 9.542651011430099
    "E20" = 13.11858250622472
    "F20" = 20.24955283636154
    "G20" = 17.06251994700407
    "H20" = 11.03694355190517
    "I20" = 15.17256391256532
    "M20" = 14.31318420190494
    "N20" = 15.7840777207902
    "O20" = 15.24150732605617
    "B21" = 12.81475104365535
    "C21" = 9.883482342676976
    "E21" = 13.31733535275646
    "F21" = 21.46857628470577
    "G21" = 17.20928646491081
    "H21" = 10.99783504210724
    "I21" = 15.06843094987005
    "M21" = 14.69219300100399
    "N21" = 15.75572313580017
    "O21" = 15.20685586694384
    "B22" = 13.21853779152592
    "C22" = 10.09899278458132
    "E22" = 13.44991976511709
    "F22" = 22.22866616901552
    "G22" = 17.31220141839842
    "H22" = 10.97415671724808
    "I22" = 15.00360546265419
    "M22" = 14.93663795230497
    "N22" = 15.73865000618901
    "O22" = 15.18821259910888
    "B23" = 13.00474171273053
    "C23" = 9.984646308915012
    "E23" = 13.37892578142741
    "F23" = 21.82633154458858
    "G23" = 17.25664416815107
    "H23" = 10.98662190482789
    "I23" = 15.03790795825662
    "M23" = 14.80649762163172
    "N23" = 15.74762902001615
    "O23" = 15.1977955275478
    "B24" = 12.15746547945982
    "C24" = 9.536961046244278
    "E24" = 13.11538264858929
    "F24" = 20.22900810905287
    "G24" = 17.06024142829225
    "H24" = 11.03761446812567
    "I24" = 15.17432037574417
    "M24" = 14.30694051772172
    "N24" = 15.78456582246786
    "O24" = 15.24214103583105
    "B25" = 11.16386244516981
    "C25" = 9.024416119748006
    "E25" = 12.84411527316082
    "F25" = 18.34778573295695
    "G25" = 16.87936605849923
    "H25" = 11.10067137947908
    "I25" = 15.33529816746625
    "M25" = 13.75691000756917
    "N25" = 15.83068821167836
    "O25" = 15.30705066398344
}

foreach ($cellAddress in $updatedValues.Keys) {
    $ws.Range($cellAddress).Value = $updatedValues[$cellAddress]
}

Write-Host "Updated $($updatedValues.Count) cells in $($ws.Name)"
